$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to text format so numeric-looking strings (e.g. "569.72")
# are not auto-converted to numbers by Excel, matching the original inline-string cells.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '63.242.03'
$ws.Range("E2").Value = '  +1.93%  '

$ws.Range("D3").Value = '2.545.73'
$ws.Range("E3").Value = '  +4.50%  '

$ws.Range("E4").Value = '  +0.11%  '

$ws.Range("D5").Value = '569.72'
$ws.Range("E5").Value = '  +2.56%  '

$ws.Range("D6").Value = '150.32'
$ws.Range("E6").Value = '  +8.41%  '

$ws.Range("E7").Value = '  +0.11%  '

$ws.Range("D8").Value = '0.588'
$ws.Range("E8").Value = '  +0.79%  '

$ws.Range("D9").Value = '2.537.78'
$ws.Range("E9").Value = '  +4.28%  '

$ws.Range("D10").Value = '0.106'
$ws.Range("E10").Value = '  +1.37%  '

$ws.Range("D11").Value = '5.73'
$ws.Range("E11").Value = '  -0.17%  '

$ws.Range("E12").Value = '  +1.19%  '

$ws.Range("D13").Value = '0.358'
$ws.Range("E13").Value = '  +2.83%  '

$ws.Range("D14").Value = '28.25'
$ws.Range("E14").Value = '  +8.36%  '

$ws.Range("D15").Value = '3.002.73'
$ws.Range("E15").Value = '  +4.65%  '

$ws.Range("D16").Value = '63.245.87'
$ws.Range("E16").Value = '  +2.06%  '

$ws.Range("D17").Value = '0.0000143'
$ws.Range("E17").Value = '  -0.22%  '

$ws.Range("D18").Value = '2.545.27'
$ws.Range("E18").Value = '  +4.54%  '

$ws.Range("D19").Value = '11.60'
$ws.Range("E19").Value = '  +4.11%  '

$ws.Range("D20").Value = '340.98'
$ws.Range("E20").Value = '  -1.27%  '

$ws.Range("D21").Value = '4.36'
$ws.Range("E21").Value = '  +3.55%  '

$ws.Range("D22").Value = '6.86'
$ws.Range("E22").Value = '  +0.45%  '

$ws.Range("E23").Value = '  +0.04%  '

$ws.Range("D24").Value = '65.96'
$ws.Range("E24").Value = '  +1.22%  '

$ws.Range("E25").Value = '  -1.52%  '

$ws.Range("D26").Value = '1.60'
$ws.Range("E26").Value = '  +5.30%  '

$ws.Range("B27").Value = 'Binance-PegBSC-USD'
$ws.Range("C27").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D27").Value = '1.00'
$ws.Range("E27").Value = '  +0.15%  '

$ws.Range("B28").Value = 'InternetComputer(DFINITY)'
$ws.Range("C28").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D28").Value = '8.43'
$ws.Range("E28").Value = '  +2.42%  '

$ws.Range("D29").Value = '1.47'
$ws.Range("E29").Value = '  +9.02%  '

$ws.Range("D30").Value = '7.10'
$ws.Range("E30").Value = '  +11.90%  '

$ws.Range("D31").Value = '0.0₃0829'
$ws.Range("E31").Value = '  +4.99%  '

$ws.Range("D32").Value = '1.86'
$ws.Range("E32").Value = '  +2.85%  '

$ws.Range("D33").Value = '176.01'
$ws.Range("E33").Value = '  +2.43%  '

$ws.Range("D34").Value = '1.55'
$ws.Range("E34").Value = '  +6.86%  '

$ws.Range("D35").Value = '418.59'
$ws.Range("E35").Value = '  +13.68%  '

$ws.Range("D36").Value = '0.406'
$ws.Range("E36").Value = '  +2.30%  '

$ws.Range("D37").Value = '19.09'
$ws.Range("E37").Value = '  +2.64%  '

$ws.Range("D38").Value = '4.42'
$ws.Range("E38").Value = '  -1.94%  '

$ws.Range("E39").Value = '  +0.02%  '

$ws.Range("D40").Value = '1.75'
$ws.Range("E40").Value = '  +2.72%  '

$ws.Range("E41").Value = '  +0.09%  '

$ws.Range("D42").Value = '40.05'
$ws.Range("E42").Value = '  +2.00%  '

$ws.Range("D43").Value = '155.05'
$ws.Range("E43").Value = '  +5.78%  '

$ws.Range("D44").Value = '3.80'
$ws.Range("E44").Value = '  +3.24%  '

$ws.Range("D45").Value = '21.04'
$ws.Range("E45").Value = '  +1.71%  '

$ws.Range("D46").Value = '0.608'
$ws.Range("E46").Value = '  +3.26%  '

$ws.Range("D47").Value = '0.0530'
$ws.Range("E47").Value = '  +2.22%  '

$ws.Range("D48").Value = '0.0965'
$ws.Range("E48").Value = '  +0.69%  '

$ws.Range("D49").Value = '0.0236'
$ws.Range("E49").Value = '  +6.40%  '

$ws.Range("D50").Value = '18.67'
$ws.Range("E50").Value = '  +4.24%  '

$ws.Range("B51").Value = 'BabyDogeCoin'
$ws.Range("C51").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D51").Value = '0.0₆0230'
$ws.Range("E51").Value = '  +6.23%  '

# Reset the style index on column D back to the workbook default (no explicit style),
# since the cells only need to retain their text type, not a visible style change.
$ws.Range("D2:D51").Style = "Normal"
